$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "Palak"
$ws.Range("D2").Value = "PL_0001"
$ws.Range("E2").Value = 20.56
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 1000

# Remove row 3 entirely (the Ginger entry) by clearing contents & formatting
# (using Clear instead of Delete avoids shifting the whole-column data
# validation ranges below the removed row)
$ws.Range("A3:H3").Clear()

# Update selection to H2 to match final state
$ws.Range("H2").Select()
